$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect the latest scrape.
# For cells whose new text looks like a plain number, the Price column
# number format is temporarily switched to Text ("@") before the value is
# written so Excel keeps the exact original text (instead of silently
# converting it to a floating point number and losing precision/zeros),
# then the format is cleared again so the cell keeps its original (default)
# style index.

$ws.Range("D2").Value = "27.680.38"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.849.51"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4261"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3631"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.73"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07297"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8745"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.58"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "1.870.91"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.314"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.506"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06891"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.67"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008993"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "27.691.71"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.967"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("D25").Value = "2.073.63"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.964"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.82"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.55"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.250"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.861"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08921"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7586"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.969"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.511"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.099"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01930"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.814"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5064"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1650"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.759"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.327"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06544"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.28"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.98"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4654"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.616"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.756"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.56%  "
